$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing Row 5 (Coritiba vs Botafogo SP) odds ---
$ws.Range("G5").Value = 1.53
$ws.Range("H5").Value = 3.7
$ws.Range("I5").Value = 7
$ws.Range("L5").Value = 7
$ws.Range("Z5").Value = 10
$ws.Range("AH5").Value = 15
$ws.Range("AI5").Value = 34
$ws.Range("AK5").Value = 81
$ws.Range("AN5").Value = 3.25
$ws.Range("AU5").Value = 10
$ws.Range("AV5").Value = 81
$ws.Range("AY5").Value = 41

# --- Update existing Row 8 (Platense Municipal vs Cacahuatique) odds ---
$ws.Range("J8").Value = 3.65
$ws.Range("L8").Value = 2.9
$ws.Range("M8").Value = 1.09
$ws.Range("U8").Value = 1.91
$ws.Range("W8").Value = 7.9
$ws.Range("X8").Value = 15
$ws.Range("AA8").Value = 30
$ws.Range("AG8").Value = 1000
$ws.Range("AO8").Value = 17
$ws.Range("AP8").Value = 27
$ws.Range("AU8").Value = 7.5
$ws.Range("AX8").Value = 4
$ws.Range("AY8").Value = 12

# --- Insert new Row 9: EL SALVADOR - PRIMERA DIVISION, Aguila vs Luis Angel Firpo ---
$ws.Rows.Item(9).Insert()
$ws.Range("A9").Value = "tKBCTArm"
$ws.Range("B9").Value = "22/11/2024"
$ws.Range("C9").Value = "22:00"
$ws.Range("D9").Value = "EL SALVADOR - PRIMERA DIVISION"
$ws.Range("E9").Value = "Aguila"
$ws.Range("F9").Value = "Luis Angel Firpo"
$ws.Range("G9").Value = 2.02
$ws.Range("H9").Value = 3.4
$ws.Range("I9").Value = 3.3
$ws.Range("J9").Value = 2.6
$ws.Range("K9").Value = 2.15
$ws.Range("L9").Value = 3.75
$ws.Range("M9").Value = 1.01
$ws.Range("N9").Value = 11
$ws.Range("O9").Value = 1.24
$ws.Range("P9").Value = 3.3
$ws.Range("Q9").Value = 1.7
$ws.Range("R9").Value = 1.91
$ws.Range("S9").Value = 1.33
$ws.Range("T9").Value = 3.04
$ws.Range("U9").Value = 1.6
$ws.Range("V9").Value = 2.07
$ws.Range("W9").Value = 8.75
$ws.Range("X9").Value = 10.75
$ws.Range("Y9").Value = 8.5
$ws.Range("Z9").Value = 19
$ws.Range("AA9").Value = 15
$ws.Range("AB9").Value = 23
$ws.Range("AC9").Value = 11.5
$ws.Range("AD9").Value = 6.7
$ws.Range("AE9").Value = 13
$ws.Range("AF9").Value = 50
$ws.Range("AG9").Value = 350
$ws.Range("AH9").Value = 11
$ws.Range("AI9").Value = 18.5
$ws.Range("AJ9").Value = 11.25
$ws.Range("AK9").Value = 45
$ws.Range("AL9").Value = 28
$ws.Range("AM9").Value = 32
$ws.Range("AN9").Value = 4
$ws.Range("AO9").Value = 10.25
$ws.Range("AP9").Value = 17
$ws.Range("AQ9").Value = 37
$ws.Range("AR9").Value = 65
$ws.Range("AS9").Value = 200
$ws.Range("AT9").Value = 2.82
$ws.Range("AU9").Value = 6.8
$ws.Range("AV9").Value = 55
$ws.Range("AW9").Value = $null
$ws.Range("AX9").Value = 5.2
$ws.Range("AY9").Value = 18
$ws.Range("AZ9").Value = 24
$ws.Range("BA9").Value = 90
$ws.Range("BB9").Value = 120
$ws.Range("BC9").Value = 300
$ws.Range("BD9").Value = $null

# --- Insert new Row 17: NEW ZEALAND - NATIONAL LEAGUE, Cashmere Technical vs Birkenhead ---
$ws.Rows.Item(17).Insert()
$ws.Range("A17").Value = "tvUXozpg"
$ws.Range("B17").Value = "22/11/2024"
$ws.Range("C17").Value = "22:00"
$ws.Range("D17").Value = "NEW ZEALAND - NATIONAL LEAGUE"
$ws.Range("E17").Value = "Cashmere Technical"
$ws.Range("F17").Value = "Birkenhead"
$ws.Range("G17").Value = 3.9
$ws.Range("H17").Value = 4.55
$ws.Range("I17").Value = 1.57
$ws.Range("J17").Value = 3.85
$ws.Range("K17").Value = 2.77
$ws.Range("L17").Value = 1.95
$ws.Range("M17").Value = $null
$ws.Range("N17").Value = $null
$ws.Range("O17").Value = 1.01
$ws.Range("P17").Value = 7.9
$ws.Range("Q17").Value = 1.21
$ws.Range("R17").Value = 3.44
$ws.Range("S17").Value = 1.14
$ws.Range("T17").Value = 4.9
$ws.Range("U17").Value = 1.28
$ws.Range("V17").Value = 3.52
$ws.Range("W17").Value = 22
$ws.Range("X17").Value = 28
$ws.Range("Y17").Value = 13
$ws.Range("Z17").Value = 50
$ws.Range("AA17").Value = 23
$ws.Range("AB17").Value = 19.5
$ws.Range("AC17").Value = 30
$ws.Range("AD17").Value = 10
$ws.Range("AE17").Value = 10.5
$ws.Range("AF17").Value = 23
$ws.Range("AG17").Value = 80
$ws.Range("AH17").Value = 14.5
$ws.Range("AI17").Value = 11.75
$ws.Range("AJ17").Value = 8.25
$ws.Range("AK17").Value = 13.5
$ws.Range("AL17").Value = 9.5
$ws.Range("AM17").Value = 12.5
$ws.Range("AN17").Value = 7.1
$ws.Range("AO17").Value = 19
$ws.Range("AP17").Value = 17
$ws.Range("AQ17").Value = 70
$ws.Range("AR17").Value = 65
$ws.Range("AS17").Value = 110
$ws.Range("AT17").Value = 5.1
$ws.Range("AU17").Value = 6.4
$ws.Range("AV17").Value = 28
$ws.Range("AW17").Value = 250
$ws.Range("AX17").Value = 4.5
$ws.Range("AY17").Value = 7.4
$ws.Range("AZ17").Value = 10.5
$ws.Range("BA17").Value = 18
$ws.Range("BB17").Value = 25
$ws.Range("BC17").Value = 70
$ws.Range("BD17").Value = $null

# --- Update Row 20 (Nacional vs Montevideo City) odds (now shifted from old row 18) ---
$ws.Range("G20").Value = 1.26
$ws.Range("H20").Value = 5.5
$ws.Range("I20").Value = 9.5
$ws.Range("K20").Value = 2.72
$ws.Range("L20").Value = 7.7
$ws.Range("N20").Value = 9.25
$ws.Range("O20").Value = 1.17
$ws.Range("P20").Value = 4.5
$ws.Range("Q20").Value = 1.52
$ws.Range("R20").Value = 2.37
$ws.Range("S20").Value = 1.25
$ws.Range("T20").Value = 3.55
$ws.Range("W20").Value = 8
$ws.Range("Z20").Value = 7.7
$ws.Range("AC20").Value = 9.25
$ws.Range("AD20").Value = 11
$ws.Range("AG20").Value = 800
$ws.Range("AH20").Value = 27
$ws.Range("AI20").Value = 70
$ws.Range("AJ20").Value = 29
$ws.Range("AK20").Value = 250
$ws.Range("AL20").Value = 110
$ws.Range("AM20").Value = 90
$ws.Range("AN20").Value = 3.25
$ws.Range("AO20").Value = 5.3
$ws.Range("AT20").Value = 3.55
$ws.Range("AV20").Value = 70
$ws.Range("AX20").Value = 10.25
$ws.Range("AY20").Value = 50
$ws.Range("AZ20").Value = 45
$ws.Range("BA20").Value = 400
$ws.Range("BB20").Value = 350
$ws.Range("BC20").Value = 500

